$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Mike's hours (column D) for Week 5 and Week 6
$ws.Range("D8").Value = 15
$ws.Range("D9").Value = 16

# Update Patrick's hours (column F) for Week 5 and Week 6
$ws.Range("F8").Value = 11
$ws.Range("F9").Value = 15

# Update the active selection to match the authored state
$ws.Range("F8").Select()
